$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells: "<col>_old" -> "<col>_FV2210", "<col>_new" -> "<col>_FV2304" ---
$leftHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
$rightHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $leftHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $leftHeaders[$i]
}
$ws.Cells.Item(1, 11).Value = "diff"
for ($i = 0; $i -lt $rightHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $rightHeaders[$i]
}

# --- 2) Turn the data range into an Excel Table ("Table1") ---
# Temporarily clear the header row's existing direct formatting before the
# table is created so Excel doesn't need to capture an overlay (dxf) for a
# header style that differs from the table's computed header look; then
# restore the original header formatting afterwards.
$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U56")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"

# Restore the header row's original look (bold, centered, wrapped, grey
# fill, thin border on all sides) now that the table exists.
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Interior.Color = 14277081
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# --- 3) Freeze the header row (split/freeze pane at row 2) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
